$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "51.593.66"
$ws.Cells.Item(2,5).Value = "  +1.53%  "
$ws.Cells.Item(3,4).Value = "3.024.58"
$ws.Cells.Item(3,5).Value = "  +3.16%  "
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "0.999"
$ws.Cells.Item(4,5).Value = "  -0.01%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "379.18"
$ws.Cells.Item(5,5).Value = "  +0.63%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "103.17"
$ws.Cells.Item(6,5).Value = "  +2.90%  "
$ws.Cells.Item(7,5).Value = "  +1.52%  "
$ws.Cells.Item(8,5).Value = "  -0.01%  "
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.596"
$ws.Cells.Item(9,5).Value = "  +3.48%  "
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "36.75"
$ws.Cells.Item(10,5).Value = "  +2.50%  "
$ws.Cells.Item(11,5).Value = "  -0.27%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.0860"
$ws.Cells.Item(12,5).Value = "  +1.01%  "
$ws.Cells.Item(13,4).Value = "3.499.26"
$ws.Cells.Item(13,5).Value = "  +2.68%  "
$ws.Cells.Item(14,5).Value = "  +1.96%  "
$ws.Cells.Item(15,5).Value = "  +2.20%  "
$ws.Cells.Item(16,4).Value = "3.034.92"
$ws.Cells.Item(16,5).Value = "  +3.26%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.984"
$ws.Cells.Item(17,5).Value = "  -0.74%  "
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "10.47"
$ws.Cells.Item(18,5).Value = "  -14.01%  "
$ws.Cells.Item(19,4).Value = "51.597.42"
$ws.Cells.Item(19,5).Value = "  +1.56%  "
$ws.Cells.Item(20,5).Value = "  +1.50%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "12.48"
$ws.Cells.Item(21,5).Value = "  +1.24%  "
$ws.Cells.Item(22,4).Value = "0.0₃0963"
$ws.Cells.Item(22,5).Value = "  +1.96%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "70.06"
$ws.Cells.Item(23,5).Value = "  +0.95%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "268.13"
$ws.Cells.Item(24,5).Value = "  +0.88%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "3.14"
$ws.Cells.Item(25,5).Value = "  -1.71%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "8.19"
$ws.Cells.Item(26,5).Value = "  +3.50%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "7.52"
$ws.Cells.Item(27,5).Value = "  +6.10%  "
$ws.Cells.Item(28,5).Value = "  +6.50%  "
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "26.22"
$ws.Cells.Item(30,5).Value = "  +2.95%  "
$ws.Cells.Item(31,5).Value = "  +1.45%  "
$ws.Cells.Item(32,5).Value = "  +3.34%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "34.34"
$ws.Cells.Item(33,5).Value = "  +3.01%  "
$ws.Cells.Item(34,2).Value = "OKB"
$ws.Cells.Item(34,3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "50.54"
$ws.Cells.Item(34,5).Value = "  +0.14%  "
$ws.Cells.Item(35,2).Value = "Toncoin"
$ws.Cells.Item(35,3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "2.05"
$ws.Cells.Item(35,5).Value = "  +0.42%  "
$ws.Cells.Item(36,2).Value = "VeChain"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.0453"
$ws.Cells.Item(36,5).Value = "  +5.35%  "
$ws.Cells.Item(37,5).Value = "  -0.09%  "
$ws.Cells.Item(38,5).Value = "  +6.52%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "17.43"
$ws.Cells.Item(39,5).Value = "  +6.09%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.288"
$ws.Cells.Item(40,5).Value = "  +11.50%  "
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "1.87"
$ws.Cells.Item(41,5).Value = "  +3.96%  "
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "2.57"
$ws.Cells.Item(42,5).Value = "  +5.59%  "
$ws.Cells.Item(43,2).Value = "Stellar"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.116"
$ws.Cells.Item(43,5).Value = "  +0.83%  "
$ws.Cells.Item(44,2).Value = "Monero"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "126.80"
$ws.Cells.Item(44,5).Value = "  +2.87%  "
$ws.Cells.Item(45,5).Value = "  +9.41%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "21.94"
$ws.Cells.Item(46,5).Value = "  +4.38%  "
$ws.Cells.Item(47,5).Value = "  +1.70%  "
$ws.Cells.Item(48,5).Value = "  +1.68%  "
$ws.Cells.Item(49,4).Value = "2.033.98"
$ws.Cells.Item(49,5).Value = "  +1.61%  "
$ws.Cells.Item(50,4).Value = "3.324.51"
$ws.Cells.Item(50,5).Value = "  +3.21%  "
$ws.Cells.Item(51,5).Value = "  +2.00%  "